$wb = $excel.ActiveWorkbook

# --- Sheet 1: Metadata ---
$ws1 = $wb.Worksheets.Item(1)

# Shift rows 11-14 down to rows 12-15 (to make room for the new "Jurisdiction" row),
# copying values first (bottom-up) and then formatting (so styles/borders carry over
# correctly without creating new style entries).
$ws1.Range("A15").Value = $ws1.Range("A14").Value2
$ws1.Range("B15").Value = $ws1.Range("B14").Value2
$ws1.Range("A14").Copy()
$ws1.Range("A15").PasteSpecial(-4122)
$ws1.Range("B14").Copy()
$ws1.Range("B15").PasteSpecial(-4122)

$ws1.Range("A14").Value = $ws1.Range("A13").Value2
$ws1.Range("B14").Value = $ws1.Range("B13").Value2
$ws1.Range("A13").Copy()
$ws1.Range("A14").PasteSpecial(-4122)
$ws1.Range("B13").Copy()
$ws1.Range("B14").PasteSpecial(-4122)

$ws1.Range("A13").Value = $ws1.Range("A12").Value2
$ws1.Range("B13").Value = $ws1.Range("B12").Value2
$ws1.Range("A12").Copy()
$ws1.Range("A13").PasteSpecial(-4122)
$ws1.Range("B12").Copy()
$ws1.Range("B13").PasteSpecial(-4122)

$ws1.Range("A12").Value = $ws1.Range("A11").Value2
$ws1.Range("B12").Value = $ws1.Range("B11").Value2
$ws1.Range("A11").Copy()
$ws1.Range("A12").PasteSpecial(-4122)
$ws1.Range("B11").Copy()
$ws1.Range("B12").PasteSpecial(-4122)

# New row 11: Jurisdiction | (blank) - reuses the style already present on A11/B11
$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""

# Update Version value (row 3) and Date value (row 8)
$ws1.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws1.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# --- Sheet 2: Include ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Include #0"
